$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target stored column widths (as they will appear in the saved <cols> xml):
#   A = 6, B = 17, C = 27, D = 17, E = 27, F = 27
#
# Excel's COM ColumnWidth setter is expressed in "characters" and the engine
# adds the standard padding (5/6 of a character) when persisting the raw
# <col width="..."> attribute, so we subtract that same padding here to land
# on the exact stored widths from the diff.
$padding = 5 / 6

$widths = @{
    1 = 6
    2 = 17
    3 = 27
    4 = 17
    5 = 27
    6 = 27
}

foreach ($col in 1..6) {
    $ws.Columns.Item($col).ColumnWidth = $widths[$col] - $padding
}
